# feat: Catch dialog and view allow the bot to gracefully handle errors.
#
# Adds two new "disambiguation" intents (disanbiguationone / disanbiguationtwo)
# plus their training utterances, and a new "error" intent, to the
# sample-complexdialogs intents sheet. Also nudges the tab-ratio of the
# workbook window and leaves the selection on the last edited cell, matching
# the author's recorded session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new "disanbiguationone" intent block (rows 38-46) ---
$ws.Range("A38").Value = "disanbiguationone"
$ws.Range("A39").Value = "je veux résilier ma seconde carte"
$ws.Range("A40").Value = "je veux résilier ma carte"
$ws.Range("A41").Value = "résiliation de ma carte"
$ws.Range("A42").Value = "comment faire une résiliation ? "
$ws.Range("A43").Value = "comment puis-je procéder?"
$ws.Range("A44").Value = "je veux accéder au formulaire"
$ws.Range("A45").Value = "je ne veux plus l'utiliser"
$ws.Range("A46").Value = "je veux arrêter "

# row 47 is left blank on purpose - it separates the two intent blocks

# --- new "disanbiguationtwo" intent block (rows 48-52) ---
$ws.Range("A48").Value = "disanbiguationtwo"
$ws.Range("A49").Value = "je veux boire"
$ws.Range("A50").Value = "je veux jouer"
$ws.Range("A51").Value = "je veux manger"
$ws.Range("A52").Value = "i see a disambiguation"

# row 53 is left blank on purpose - it separates the two intent blocks

# --- new "error" intent block (rows 54-55) ---
$ws.Range("A54").Value = "error"
$ws.Range("A55").Value = "error"

# Match the author's recorded window/view state.
$excel.ActiveWindow.TabRatio = 993

# Leave the selection where the author left it after typing the last entry.
[void]$ws.Range("G49").Select()
